$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos data: set each changed cell as text
# (prefixed with an apostrophe to stop Excel re-typing numeric-looking
# strings as Number), then reset the cell style to Normal so no stray
# formatting (e.g. quote-prefix / text number-format) is introduced.

$ws.Range("D2").Value = "'27.221.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.29%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.895.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.04%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'306.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.34%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.03%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.5393"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.66%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.76%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07276"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.12%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'21.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.80%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.8973"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.35%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08184"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.32%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'94.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.52%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.36%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'1.721.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -11.83%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +1.64%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008634"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.18%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.07%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'27.024.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.55%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.033"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.15%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.79%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.463"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.60%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'148.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.07%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.296"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.63%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'18.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.88%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.752"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.52%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'116.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.71%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'4.820"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.28%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.643"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.43%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.09164"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.52%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.8196"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.00%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05041"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.28%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.04%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.020"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.26%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -4.03%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.670"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.31%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.5959"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.13%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01983"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.89%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.074"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.41%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'9.224"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.27%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'6.629"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.02%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'114.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.02%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.5093"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.64%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.1525"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.64%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'PaxDollar"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.09%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'10.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.33%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.628"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.44%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'37.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.13%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06086"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.67%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'62.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.68%  "
$ws.Range("E51").Style = "Normal"
